$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 229.25
$ws.Range("I15").Value = 229.25
$ws.Range("K15").Value = 687.75
$ws.Range("M15").Value = -518.75
$ws.Range("H40").Value = 8198688.5
$ws.Range("I40").Value = 2009.5962
$ws.Range("J40").Value = 55557276
$ws.Range("K40").Value = 2009.5962
$ws.Range("L40").Value = 55557276
$ws.Range("M40").Value = -1834.5962
$ws.Range("N40").Value = -55557626
$ws.Range("H41").Value = 457.2857
$ws.Range("I41").Value = 450.2
$ws.Range("J41").Value = 475
$ws.Range("K41").Value = 450.2
$ws.Range("L41").Value = 475
$ws.Range("M41").Value = -10.19999999999999
$ws.Range("N41").Value = -1355
$ws.Range("H55").Value = 461.33334
$ws.Range("I55").Value = 815
$ws.Range("J55").Value = 178.4
$ws.Range("K55").Value = 815
$ws.Range("L55").Value = 178.4
$ws.Range("M55").Value = -601
$ws.Range("N55").Value = -606.4
$ws.Range("H64").Value = 6254.6665
$ws.Range("I64").Value = 7582
$ws.Range("J64").Value = 3600
$ws.Range("K64").Value = 7582
$ws.Range("L64").Value = 3600
$ws.Range("M64").Value = -7334
$ws.Range("N64").Value = -4096
$ws.Range("H67").Value = 6254.6665
$ws.Range("I67").Value = 7582
$ws.Range("J67").Value = 3600
$ws.Range("K67").Value = 7582
$ws.Range("L67").Value = 3600
$ws.Range("M67").Value = -6724
$ws.Range("N67").Value = -5316
$ws.Range("H100").Value = 41667920
$ws.Range("I100").Value = 41667920
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 41667920
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -41667379
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 22728372
$ws.Range("I107").Value = 25000708
$ws.Range("K107").Value = 25000708
$ws.Range("M107").Value = -24998788
$ws.Range("H115").Value = 306.66666
$ws.Range("I115").Value = 306.66666
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 919.9999799999999
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 647.0000200000001
$ws.Range("N115").ClearContents()
$ws.Range("H132").Value = 2117.975
$ws.Range("I132").Value = 1498.5588
$ws.Range("J132").Value = 5628
$ws.Range("K132").Value = 4495.6764
$ws.Range("L132").Value = 16884
$ws.Range("M132").Value = -1965.6764
$ws.Range("N132").Value = -21944
$ws.Range("H137").Value = 1939.0344
$ws.Range("I137").Value = 1378.963
$ws.Range("J137").Value = 9500
$ws.Range("K137").Value = 4136.889
$ws.Range("L137").Value = 28500
$ws.Range("M137").Value = -1586.889
$ws.Range("N137").Value = -33600
$ws.Range("H138").Value = 2836.4348
$ws.Range("I138").Value = 1454.1852
$ws.Range("J138").Value = 3725.024
$ws.Range("K138").Value = 4362.5556
$ws.Range("L138").Value = 11175.072
$ws.Range("M138").Value = 777.4444000000003
$ws.Range("N138").Value = -21455.072

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8766.947
$ws.Range("I32").Value = 7460.3213
$ws.Range("J32").Value = 19742.6
$ws.Range("K32").Value = 7460.3213
$ws.Range("L32").Value = 19742.6
$ws.Range("M32").Value = -7173.3213
$ws.Range("N32").Value = -20316.6

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 100002610
$ws.Range("I99").Value = 500000600
$ws.Range("J99").Value = 3110.875
$ws.Range("K99").Value = 500000600
$ws.Range("L99").Value = 3110.875
$ws.Range("M99").Value = -499999102
$ws.Range("N99").Value = -6106.875
$ws.Range("H134").Value = 27935.582
$ws.Range("I134").Value = 5024.919
$ws.Range("K134").Value = 15074.757
$ws.Range("M134").Value = -12539.757

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 260695.06
$ws.Range("I31").Value = 1207.7675
$ws.Range("J31").Value = 917045.3
$ws.Range("K31").Value = 1207.7675
$ws.Range("L31").Value = 917045.3
$ws.Range("M31").Value = -912.7674999999999
$ws.Range("N31").Value = -917635.3
$ws.Range("H34").Value = 260695.06
$ws.Range("I34").Value = 1207.7675
$ws.Range("J34").Value = 917045.3
$ws.Range("K34").Value = 1207.7675
$ws.Range("L34").Value = 917045.3
$ws.Range("M34").Value = -1005.7675
$ws.Range("N34").Value = -917449.3
$ws.Range("H132").Value = 2605.3
$ws.Range("I132").Value = 1863.6316
$ws.Range("J132").Value = 3886.3635
$ws.Range("K132").Value = 5590.8948
$ws.Range("L132").Value = 11659.0905
$ws.Range("M132").Value = -3060.8948
$ws.Range("N132").Value = -16719.0905

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 57.75
$ws.Range("I14").Value = 57.75
$ws.Range("K14").Value = 173.25
$ws.Range("M14").Value = -0.25
$ws.Range("H92").Value = 773.3684
$ws.Range("I92").Value = 668.1667
$ws.Range("J92").Value = 821.9231
$ws.Range("K92").Value = 2004.5001
$ws.Range("L92").Value = 2465.7693
$ws.Range("M92").Value = -756.5001
$ws.Range("N92").Value = -4961.7693
$ws.Range("H113").Value = 3750510.8
$ws.Range("I113").Value = 7143283.5
$ws.Range("J113").Value = 1111687.5
$ws.Range("K113").Value = 21429850.5
$ws.Range("L113").Value = 3335062.5
$ws.Range("M113").Value = -21427680.5
$ws.Range("N113").Value = -3339402.5
$ws.Range("H132").Value = 1034139.2
$ws.Range("J132").Value = 1636625.2
$ws.Range("L132").Value = 14729626.8
$ws.Range("N132").Value = -14734686.8

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1279.1333
$ws.Range("I107").Value = 355.5
$ws.Range("J107").Value = 2334.7144
$ws.Range("K107").Value = 355.5
$ws.Range("L107").Value = 2334.7144
$ws.Range("M107").Value = 1564.5
$ws.Range("N107").Value = -6174.7144
$ws.Range("H132").Value = 5159
$ws.Range("I132").Value = 5647.6
$ws.Range("J132").Value = 4635.5
$ws.Range("K132").Value = 16942.8
$ws.Range("L132").Value = 13906.5
$ws.Range("M132").Value = -14412.8
$ws.Range("N132").Value = -18966.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 731.3333
$ws.Range("I46").Value = 462.66666
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 462.66666
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -274.66666
$ws.Range("N46").Value = -1376
$ws.Range("H55").Value = 152
$ws.Range("I55").Value = 193.33333
$ws.Range("J55").Value = 90
$ws.Range("K55").Value = 193.33333
$ws.Range("L55").Value = 90
$ws.Range("M55").Value = -20.33332999999999
$ws.Range("N55").Value = -436
$ws.Range("H122").Value = 9055930
$ws.Range("I122").Value = 8937921
$ws.Range("J122").Value = 10000000
$ws.Range("K122").Value = 26813763
$ws.Range("L122").Value = 30000000
$ws.Range("M122").Value = -26811313
$ws.Range("N122").Value = -30004900
$ws.Range("H132").Value = 7412651.5
$ws.Range("I132").Value = 7941865
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 23825595
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -23823065
$ws.Range("N132").Value = -16058.9999
$ws.Range("H136").Value = 10102.581
$ws.Range("I136").Value = 6660.8076
$ws.Range("J136").Value = 27999.8
$ws.Range("K136").Value = 19982.4228
$ws.Range("L136").Value = 83999.39999999999
$ws.Range("M136").Value = -17432.4228
$ws.Range("N136").Value = -89099.39999999999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10000000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H107").Value = 111111850
$ws.Range("I107").Value = 200000530
$ws.Range("K107").Value = 600001590
$ws.Range("M107").Value = -599999670
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 3500
$ws.Range("K113").Value = 10500
$ws.Range("M113").Value = -8330
$ws.Range("H132").Value = 2346.3872
$ws.Range("I132").Value = 1635.238
$ws.Range("J132").Value = 3839.8
$ws.Range("K132").Value = 4905.714
$ws.Range("L132").Value = 11519.4
$ws.Range("M132").Value = -2375.714
$ws.Range("N132").Value = -16579.4
$ws.Range("H136").Value = 2710.635
$ws.Range("I136").Value = 2700.8064
$ws.Range("K136").Value = 8102.4192
$ws.Range("M136").Value = -5552.4192
